$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Edificações - Campus Congonhas): Inscritos, Pagos, Inscrições homologadas
$ws.Range("E6").Value = 48
$ws.Range("F6").Value = 23
$ws.Range("H6").Value = 23

# Row 9 (Eletrotécnica - Campus Conselheiro Lafaiete): Inscritos
$ws.Range("E9").Value = 17

# Row 12 (Metalurgia - Campus Ouro Branco): Inscritos
$ws.Range("E12").Value = 22

# Row 15 (Campus Ouro Preto group): Pagos, Inscrições homologadas
$ws.Range("F15").Value = 40
$ws.Range("H15").Value = 40

# Row 16 (Campus Ouro Preto group): Inscritos, Pagos, Inscrições homologadas
$ws.Range("E16").Value = 281
$ws.Range("F16").Value = 79
$ws.Range("H16").Value = 79
